{"js": "// Replace the date and each \"NNN\u00d7N=\" multiplication-drill cell with its\n// updated value. Old values are unique in the document, so a literal,\n// case-sensitive search-and-replace on each pair is safe and unambiguous.\nconst replacements = [\n  [\"2024-06-20 Thursday\", \"2024-06-21 Friday\"],\n  [\"498\u00d75=\", \"355\u00d78=\"],\n  [\"135\u00d78=\", \"673\u00d77=\"],\n  [\"890\u00d77=\", \"573\u00d76=\"],\n  [\"418\u00d79=\", \"230\u00d72=\"],\n  [\"385\u00d79=\", \"638\u00d78=\"],\n  [\"373\u00d77=\", \"216\u00d74=\"],\n  [\"685\u00d74=\", \"599\u00d73=\"],\n  [\"238\u00d79=\", \"824\u00d72=\"],\n  [\"987\u00d75=\", \"787\u00d78=\"],\n  [\"491\u00d75=\", \"815\u00d74=\"],\n  [\"294\u00d76=\", \"398\u00d72=\"],\n  [\"334\u00d74=\", \"271\u00d77=\"],\n  [\"297\u00d72=\", \"943\u00d79=\"],\n  [\"545\u00d75=\", \"762\u00d76=\"],\n  [\"109\u00d73=\", \"406\u00d77=\"],\n  [\"287\u00d77=\", \"494\u00d78=\"],\n  [\"210\u00d73=\", \"924\u00d79=\"],\n  [\"892\u00d79=\", \"307\u00d79=\"],\n  [\"598\u00d72=\", \"239\u00d72=\"],\n  [\"413\u00d79=\", \"846\u00d72=\"],\n  [\"895\u00d72=\", \"589\u00d72=\"],\n  [\"309\u00d79=\", \"740\u00d79=\"],\n  [\"365\u00d76=\", \"460\u00d74=\"],\n  [\"544\u00d78=\", \"439\u00d74=\"],\n  [\"508\u00d76=\", \"391\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each \"NNN\u00d7N=\" multiplication-drill cell with its\n# updated value. Old values are unique in the document, so a literal,\n# case-sensitive Find/Replace on each pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-06-20 Thursday\", \"2024-06-21 Friday\"),\n  @(\"498\u00d75=\", \"355\u00d78=\"),\n  @(\"135\u00d78=\", \"673\u00d77=\"),\n  @(\"890\u00d77=\", \"573\u00d76=\"),\n  @(\"418\u00d79=\", \"230\u00d72=\"),\n  @(\"385\u00d79=\", \"638\u00d78=\"),\n  @(\"373\u00d77=\", \"216\u00d74=\"),\n  @(\"685\u00d74=\", \"599\u00d73=\"),\n  @(\"238\u00d79=\", \"824\u00d72=\"),\n  @(\"987\u00d75=\", \"787\u00d78=\"),\n  @(\"491\u00d75=\", \"815\u00d74=\"),\n  @(\"294\u00d76=\", \"398\u00d72=\"),\n  @(\"334\u00d74=\", \"271\u00d77=\"),\n  @(\"297\u00d72=\", \"943\u00d79=\"),\n  @(\"545\u00d75=\", \"762\u00d76=\"),\n  @(\"109\u00d73=\", \"406\u00d77=\"),\n  @(\"287\u00d77=\", \"494\u00d78=\"),\n  @(\"210\u00d73=\", \"924\u00d79=\"),\n  @(\"892\u00d79=\", \"307\u00d79=\"),\n  @(\"598\u00d72=\", \"239\u00d72=\"),\n  @(\"413\u00d79=\", \"846\u00d72=\"),\n  @(\"895\u00d72=\", \"589\u00d72=\"),\n  @(\"309\u00d79=\", \"740\u00d79=\"),\n  @(\"365\u00d76=\", \"460\u00d74=\"),\n  @(\"544\u00d78=\", \"439\u00d74=\"),\n  @(\"508\u00d76=\", \"391\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n  # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
